$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 387.33334
$ws.Cells.Item(5, 9).Value = 364.9
$ws.Cells.Item(5, 10).Value = 499.5
$ws.Cells.Item(5, 11).Value = 364.9
$ws.Cells.Item(5, 12).Value = 499.5
$ws.Cells.Item(5, 13).Value = -249.9
$ws.Cells.Item(5, 14).Value = -729.5
$ws.Cells.Item(12, 8).Value = 96.818184
$ws.Cells.Item(12, 10).Value = 87.5
$ws.Cells.Item(12, 12).Value = 87.5
$ws.Cells.Item(12, 14).Value = -427.5
$ws.Cells.Item(15, 8).Value = 672.08
$ws.Cells.Item(15, 9).Value = 672.08
$ws.Cells.Item(15, 11).Value = 2016.24
$ws.Cells.Item(15, 13).Value = -1847.24
$ws.Cells.Item(40, 8).Value = 4854.476
$ws.Cells.Item(40, 9).Value = 9138.444
$ws.Cells.Item(40, 10).Value = 1641.5
$ws.Cells.Item(40, 11).Value = 9138.444
$ws.Cells.Item(40, 12).Value = 1641.5
$ws.Cells.Item(40, 13).Value = -8963.444
$ws.Cells.Item(40, 14).Value = -1991.5
$ws.Cells.Item(41, 8).Value = 1077.3478
$ws.Cells.Item(41, 10).Value = 419
$ws.Cells.Item(41, 12).Value = 419
$ws.Cells.Item(41, 14).Value = -1299
$ws.Cells.Item(43, 8).Value = 26715.385
$ws.Cells.Item(43, 9).Value = 29150
$ws.Cells.Item(43, 10).Value = 24628.572
$ws.Cells.Item(43, 11).Value = 29150
$ws.Cells.Item(43, 12).Value = 24628.572
$ws.Cells.Item(43, 13).Value = -29081
$ws.Cells.Item(43, 14).Value = -24766.572
$ws.Cells.Item(64, 8).Value = 5256.6
$ws.Cells.Item(64, 9).Value = 3402.5264
$ws.Cells.Item(64, 10).Value = 8459.091
$ws.Cells.Item(64, 11).Value = 3402.5264
$ws.Cells.Item(64, 12).Value = 8459.091
$ws.Cells.Item(64, 13).Value = -3154.5264
$ws.Cells.Item(64, 14).Value = -8955.091
$ws.Cells.Item(67, 8).Value = 5256.6
$ws.Cells.Item(67, 9).Value = 3402.5264
$ws.Cells.Item(67, 10).Value = 8459.091
$ws.Cells.Item(67, 11).Value = 3402.5264
$ws.Cells.Item(67, 12).Value = 8459.091
$ws.Cells.Item(67, 13).Value = -2544.5264
$ws.Cells.Item(67, 14).Value = -10175.091
$ws.Cells.Item(70, 8).Value = 76571.36
$ws.Cells.Item(70, 9).Value = 8559.8
$ws.Cells.Item(70, 11).Value = 25679.4
$ws.Cells.Item(70, 13).Value = -25409.4
$ws.Cells.Item(73, 8).Value = 76571.36
$ws.Cells.Item(73, 9).Value = 8559.8
$ws.Cells.Item(73, 11).Value = 25679.4
$ws.Cells.Item(73, 13).Value = -24743.4
$ws.Cells.Item(86, 8).Value = 2466.5
$ws.Cells.Item(86, 9).Value = 1623.6471
$ws.Cells.Item(86, 10).Value = 4513.4287
$ws.Cells.Item(86, 11).Value = 1623.6471
$ws.Cells.Item(86, 12).Value = 4513.4287
$ws.Cells.Item(86, 13).Value = -500.6470999999999
$ws.Cells.Item(86, 14).Value = -6759.4287
$ws.Cells.Item(88, 8).Value = 3633
$ws.Cells.Item(88, 9).Value = 2949.5
$ws.Cells.Item(88, 11).Value = 2949.5
$ws.Cells.Item(88, 13).Value = -2543.5
$ws.Cells.Item(89, 8).Value = 2466.5
$ws.Cells.Item(89, 9).Value = 1623.6471
$ws.Cells.Item(89, 10).Value = 4513.4287
$ws.Cells.Item(89, 11).Value = 8118.2355
$ws.Cells.Item(89, 12).Value = 22567.1435
$ws.Cells.Item(89, 13).Value = -2502.2355
$ws.Cells.Item(89, 14).Value = -33799.14350000001
$ws.Cells.Item(91, 8).Value = 3633
$ws.Cells.Item(91, 9).Value = 2949.5
$ws.Cells.Item(91, 11).Value = 2949.5
$ws.Cells.Item(91, 13).Value = -1545.5
$ws.Cells.Item(98, 8).Value = 3310.5454
$ws.Cells.Item(98, 10).Value = 7637
$ws.Cells.Item(98, 12).Value = 7637
$ws.Cells.Item(98, 14).Value = -10633
$ws.Cells.Item(106, 8).Value = 1284.1666
$ws.Cells.Item(106, 9).Value = 1284.1666
$ws.Cells.Item(106, 11).Value = 1284.1666
$ws.Cells.Item(106, 13).Value = -653.1666
$ws.Cells.Item(116, 8).Value = 6117
$ws.Cells.Item(116, 9).Value = 4803.7856
$ws.Cells.Item(116, 10).Value = 9181.167
$ws.Cells.Item(116, 11).Value = 4803.7856
$ws.Cells.Item(116, 12).Value = 9181.167
$ws.Cells.Item(116, 13).Value = -1361.7856
$ws.Cells.Item(116, 14).Value = -16065.167
$ws.Cells.Item(122, 8).Value = 3310.5454
$ws.Cells.Item(122, 10).Value = 7637
$ws.Cells.Item(122, 12).Value = 22911
$ws.Cells.Item(122, 14).Value = -27811
$ws.Cells.Item(129, 8).Value = 50001644
$ws.Cells.Item(129, 9).Value = 66668452
$ws.Cells.Item(129, 11).Value = 200005356
$ws.Cells.Item(129, 13).Value = -200000356
$ws.Cells.Item(131, 8).Value = 16765.416
$ws.Cells.Item(131, 9).Value = 3227.4285
$ws.Cells.Item(131, 11).Value = 9682.2855
$ws.Cells.Item(131, 13).Value = -4642.2855
$ws.Cells.Item(135, 8).Value = 1772.8462
$ws.Cells.Item(135, 9).Value = 977.1053
$ws.Cells.Item(135, 10).Value = 3932.7144
$ws.Cells.Item(135, 11).Value = 8793.9477
$ws.Cells.Item(135, 12).Value = 35394.4296
$ws.Cells.Item(135, 13).Value = -6258.947700000001
$ws.Cells.Item(135, 14).Value = -40464.4296
$ws.Cells.Item(137, 8).Value = 1972.4375
$ws.Cells.Item(137, 9).Value = 1914.7273
$ws.Cells.Item(137, 10).Value = 2099.4
$ws.Cells.Item(137, 11).Value = 5744.1819
$ws.Cells.Item(137, 12).Value = 6298.200000000001
$ws.Cells.Item(137, 13).Value = -3194.1819
$ws.Cells.Item(137, 14).Value = -11398.2
$ws.Cells.Item(138, 8).Value = 6253010
$ws.Cells.Item(138, 9).Value = 1531.7142
$ws.Cells.Item(138, 10).Value = 7579081
$ws.Cells.Item(138, 11).Value = 4595.142599999999
$ws.Cells.Item(138, 12).Value = 22737243
$ws.Cells.Item(138, 13).Value = 544.8574000000008
$ws.Cells.Item(138, 14).Value = -22747523
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 5513.6665
$ws.Cells.Item(45, 9).Value = 6821.778
$ws.Cells.Item(45, 10).Value = 1589.3334
$ws.Cells.Item(45, 11).Value = 6821.778
$ws.Cells.Item(45, 12).Value = 1589.3334
$ws.Cells.Item(45, 13).Value = -6444.778
$ws.Cells.Item(45, 14).Value = -2343.3334
$ws.Cells.Item(61, 8).Value = 5211.2354
$ws.Cells.Item(61, 9).Value = 4289.273
$ws.Cells.Item(61, 11).Value = 4289.273
$ws.Cells.Item(61, 13).Value = -4077.273
$ws.Cells.Item(69, 8).Value = 497995
$ws.Cells.Item(69, 10).Value = 497995
$ws.Cells.Item(69, 12).Value = 497995
$ws.Cells.Item(69, 14).Value = -499493
$ws.Cells.Item(72, 8).Value = 497995
$ws.Cells.Item(72, 10).Value = 497995
$ws.Cells.Item(72, 12).Value = 1493985
$ws.Cells.Item(72, 14).Value = -1501473
$ws.Cells.Item(74, 8).Value = 57072.43
$ws.Cells.Item(74, 9).Value = 78736.766
$ws.Cells.Item(74, 11).Value = 78736.766
$ws.Cells.Item(74, 13).Value = -77862.766
$ws.Cells.Item(77, 8).Value = 57072.43
$ws.Cells.Item(77, 9).Value = 78736.766
$ws.Cells.Item(77, 11).Value = 393683.83
$ws.Cells.Item(77, 13).Value = -389315.83
$ws.Cells.Item(104, 8).Value = 85409.43
$ws.Cells.Item(104, 10).Value = 85409.43
$ws.Cells.Item(104, 12).Value = 85409.43
$ws.Cells.Item(104, 14).Value = -92397.43
$ws.Cells.Item(125, 8).Value = 29499.125
$ws.Cells.Item(125, 10).Value = 29499.125
$ws.Cells.Item(125, 12).Value = 29499.125
$ws.Cells.Item(125, 14).Value = -39339.125
$ws.Cells.Item(131, 8).Value = 164999
$ws.Cells.Item(131, 10).Value = 164999
$ws.Cells.Item(131, 12).Value = 164999
$ws.Cells.Item(131, 14).Value = -175079
$ws.Cells.Item(132, 8).Value = 4114.1313
$ws.Cells.Item(132, 9).Value = 3917.3794
$ws.Cells.Item(132, 10).Value = 4748.1113
$ws.Cells.Item(132, 11).Value = 11752.1382
$ws.Cells.Item(132, 12).Value = 14244.3339
$ws.Cells.Item(132, 13).Value = -9222.1382
$ws.Cells.Item(132, 14).Value = -19304.3339
$ws.Cells.Item(135, 8).Value = 61999.555
$ws.Cells.Item(135, 10).Value = 61999.555
$ws.Cells.Item(135, 12).Value = 61999.555
$ws.Cells.Item(135, 14).Value = -72139.555
$ws.Cells.Item(136, 8).Value = 5211.2354
$ws.Cells.Item(136, 9).Value = 4289.273
$ws.Cells.Item(136, 11).Value = 12867.819
$ws.Cells.Item(136, 13).Value = -10317.819
$ws.Cells.Item(139, 8).Value = 64399.8
$ws.Cells.Item(139, 10).Value = 64399.8
$ws.Cells.Item(139, 12).Value = 64399.8
$ws.Cells.Item(139, 14).Value = -74679.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 1295.6
$ws.Cells.Item(7, 10).Value = 396
$ws.Cells.Item(7, 12).Value = 396
$ws.Cells.Item(7, 14).Value = -622
$ws.Cells.Item(20, 8).Value = 5295.579
$ws.Cells.Item(20, 9).Value = 4197.615
$ws.Cells.Item(20, 10).Value = 7674.5
$ws.Cells.Item(20, 11).Value = 4197.615
$ws.Cells.Item(20, 12).Value = 7674.5
$ws.Cells.Item(20, 13).Value = -3950.615
$ws.Cells.Item(20, 14).Value = -8168.5
$ws.Cells.Item(22, 8).Value = 337940.47
$ws.Cells.Item(22, 10).Value = 490707.25
$ws.Cells.Item(22, 12).Value = 490707.25
$ws.Cells.Item(22, 14).Value = -491053.25
$ws.Cells.Item(86, 8).Value = 3039.6177
$ws.Cells.Item(86, 9).Value = 3710.318
$ws.Cells.Item(86, 10).Value = 1810
$ws.Cells.Item(86, 11).Value = 3710.318
$ws.Cells.Item(86, 12).Value = 1810
$ws.Cells.Item(86, 13).Value = -2587.318
$ws.Cells.Item(86, 14).Value = -4056
$ws.Cells.Item(89, 8).Value = 3039.6177
$ws.Cells.Item(89, 9).Value = 3710.318
$ws.Cells.Item(89, 10).Value = 1810
$ws.Cells.Item(89, 11).Value = 18551.59
$ws.Cells.Item(89, 12).Value = 9050
$ws.Cells.Item(89, 13).Value = -12935.59
$ws.Cells.Item(89, 14).Value = -20282
$ws.Cells.Item(103, 8).Value = 13033.333
$ws.Cells.Item(103, 10).Value = 13033.333
$ws.Cells.Item(103, 12).Value = 13033.333
$ws.Cells.Item(103, 14).Value = -15377.333
$ws.Cells.Item(105, 8).Value = 2376.6428
$ws.Cells.Item(105, 9).Value = 2378.4707
$ws.Cells.Item(105, 10).Value = 2368.875
$ws.Cells.Item(105, 11).Value = 2378.4707
$ws.Cells.Item(105, 12).Value = 2368.875
$ws.Cells.Item(105, 13).Value = -631.4706999999999
$ws.Cells.Item(105, 14).Value = -5862.875
$ws.Cells.Item(106, 8).Value = 2585.5
$ws.Cells.Item(106, 10).Value = 2585.5
$ws.Cells.Item(106, 12).Value = 2585.5
$ws.Cells.Item(106, 14).Value = -5109.5
$ws.Cells.Item(134, 8).Value = 3626.1
$ws.Cells.Item(134, 9).Value = 2323
$ws.Cells.Item(134, 10).Value = 6666.6665
$ws.Cells.Item(134, 11).Value = 6969
$ws.Cells.Item(134, 12).Value = 19999.9995
$ws.Cells.Item(134, 13).Value = -4434
$ws.Cells.Item(134, 14).Value = -25069.9995
$ws.Cells.Item(135, 8).Value = 109125
$ws.Cells.Item(135, 10).Value = 109125
$ws.Cells.Item(135, 12).Value = 109125
$ws.Cells.Item(135, 14).Value = -119265
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 325.25
$ws.Cells.Item(7, 10).Value = 556
$ws.Cells.Item(7, 12).Value = 556
$ws.Cells.Item(7, 14).Value = -782
$ws.Cells.Item(22, 8).Value = 887.1667
$ws.Cells.Item(22, 9).Value = 299.75
$ws.Cells.Item(22, 11).Value = 299.75
$ws.Cells.Item(22, 13).Value = 50.25
$ws.Cells.Item(31, 8).Value = 66907.5
$ws.Cells.Item(31, 9).Value = 94798.63
$ws.Cells.Item(31, 10).Value = 5547
$ws.Cells.Item(31, 11).Value = 94798.63
$ws.Cells.Item(31, 12).Value = 5547
$ws.Cells.Item(31, 13).Value = -94503.63
$ws.Cells.Item(31, 14).Value = -6137
$ws.Cells.Item(34, 8).Value = 66907.5
$ws.Cells.Item(34, 9).Value = 94798.63
$ws.Cells.Item(34, 10).Value = 5547
$ws.Cells.Item(34, 11).Value = 94798.63
$ws.Cells.Item(34, 12).Value = 5547
$ws.Cells.Item(34, 13).Value = -94596.63
$ws.Cells.Item(34, 14).Value = -5951
$ws.Cells.Item(58, 8).Value = 3100.1667
$ws.Cells.Item(58, 9).Value = 2740.5
$ws.Cells.Item(58, 10).Value = 3280
$ws.Cells.Item(58, 11).Value = 2740.5
$ws.Cells.Item(58, 12).Value = 3280
$ws.Cells.Item(58, 13).Value = -2537.5
$ws.Cells.Item(58, 14).Value = -3686
$ws.Cells.Item(103, 8).Value = 2302.2
$ws.Cells.Item(103, 9).Value = 2302.2
$ws.Cells.Item(103, 11).Value = 2302.2
$ws.Cells.Item(103, 13).Value = -1130.2
$ws.Cells.Item(122, 8).Value = 2105.4119
$ws.Cells.Item(122, 9).Value = 1149.4166
$ws.Cells.Item(122, 10).Value = 4399.8
$ws.Cells.Item(122, 11).Value = 3448.2498
$ws.Cells.Item(122, 12).Value = 13199.4
$ws.Cells.Item(122, 13).Value = -998.2498
$ws.Cells.Item(122, 14).Value = -18099.4
$ws.Cells.Item(124, 8).Value = 48150
$ws.Cells.Item(124, 10).Value = 48150
$ws.Cells.Item(124, 12).Value = 48150
$ws.Cells.Item(124, 14).Value = -53060
$ws.Cells.Item(125, 8).Value = 56770.832
$ws.Cells.Item(125, 10).Value = 56770.832
$ws.Cells.Item(125, 12).Value = 56770.832
$ws.Cells.Item(125, 14).Value = -61690.832
$ws.Cells.Item(132, 8).Value = 4070.55
$ws.Cells.Item(132, 9).Value = 3223.647
$ws.Cells.Item(132, 11).Value = 9670.940999999999
$ws.Cells.Item(132, 13).Value = -7140.940999999999
$ws.Cells.Item(134, 8).Value = 13079.131
$ws.Cells.Item(134, 9).Value = 8819.066
$ws.Cells.Item(134, 11).Value = 26457.198
$ws.Cells.Item(134, 13).Value = -23922.198
$ws.Cells.Item(136, 8).Value = 3100.1667
$ws.Cells.Item(136, 9).Value = 2740.5
$ws.Cells.Item(136, 10).Value = 3280
$ws.Cells.Item(136, 11).Value = 8221.5
$ws.Cells.Item(136, 12).Value = 9840
$ws.Cells.Item(136, 13).Value = -5671.5
$ws.Cells.Item(136, 14).Value = -14940
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 7441.75
$ws.Cells.Item(3, 9).Value = 2886
$ws.Cells.Item(3, 11).Value = 8658
$ws.Cells.Item(3, 13).Value = -8546
$ws.Cells.Item(5, 8).Value = 826.6429
$ws.Cells.Item(5, 9).Value = 826.6429
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 2479.9287
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -2367.9287
$ws.Cells.Item(5, 14).Value = ""
$ws.Cells.Item(18, 8).Value = 904.6
$ws.Cells.Item(18, 9).Value = 872.5
$ws.Cells.Item(18, 11).Value = 2617.5
$ws.Cells.Item(18, 13).Value = -2448.5
$ws.Cells.Item(26, 8).Value = 1093.5652
$ws.Cells.Item(26, 9).Value = 479.6
$ws.Cells.Item(26, 10).Value = 1264.1111
$ws.Cells.Item(26, 11).Value = 1438.8
$ws.Cells.Item(26, 12).Value = 3792.3333
$ws.Cells.Item(26, 13).Value = -1150.8
$ws.Cells.Item(26, 14).Value = -4368.3333
$ws.Cells.Item(29, 8).Value = 743.3333
$ws.Cells.Item(29, 9).Value = 180
$ws.Cells.Item(29, 10).Value = 1025
$ws.Cells.Item(29, 11).Value = 540
$ws.Cells.Item(29, 12).Value = 3075
$ws.Cells.Item(29, 13).Value = -263
$ws.Cells.Item(29, 14).Value = -3629
$ws.Cells.Item(32, 8).Value = 15000
$ws.Cells.Item(32, 10).Value = 15000
$ws.Cells.Item(32, 12).Value = 45000
$ws.Cells.Item(32, 14).Value = -45566
$ws.Cells.Item(44, 8).Value = 2777.4
$ws.Cells.Item(44, 9).Value = 1295.8334
$ws.Cells.Item(44, 10).Value = 4999.75
$ws.Cells.Item(44, 11).Value = 3887.5002
$ws.Cells.Item(44, 12).Value = 14999.25
$ws.Cells.Item(44, 13).Value = -3489.5002
$ws.Cells.Item(44, 14).Value = -15795.25
$ws.Cells.Item(45, 8).Value = 12502054
$ws.Cells.Item(45, 10).Value = 2538
$ws.Cells.Item(45, 12).Value = 7614
$ws.Cells.Item(45, 14).Value = -8678
$ws.Cells.Item(49, 8).Value = 1281.6666
$ws.Cells.Item(49, 9).Value = 845
$ws.Cells.Item(49, 10).Value = 1500
$ws.Cells.Item(49, 11).Value = 2535
$ws.Cells.Item(49, 12).Value = 4500
$ws.Cells.Item(49, 13).Value = -2379
$ws.Cells.Item(49, 14).Value = -4812
$ws.Cells.Item(102, 8).Value = 3875
$ws.Cells.Item(102, 9).Value = 1993
$ws.Cells.Item(102, 10).Value = 4502.3335
$ws.Cells.Item(102, 11).Value = 5979
$ws.Cells.Item(102, 12).Value = 13507.0005
$ws.Cells.Item(102, 13).Value = -3545
$ws.Cells.Item(102, 14).Value = -18375.0005
$ws.Cells.Item(108, 8).Value = 2745.375
$ws.Cells.Item(108, 9).Value = 1294.8572
$ws.Cells.Item(108, 10).Value = 12899
$ws.Cells.Item(108, 11).Value = 3884.5716
$ws.Cells.Item(108, 12).Value = 38697
$ws.Cells.Item(108, 13).Value = -1004.5716
$ws.Cells.Item(108, 14).Value = -44457
$ws.Cells.Item(110, 8).Value = 24174.334
$ws.Cells.Item(110, 9).Value = 24174.334
$ws.Cells.Item(110, 11).Value = 72523.002
$ws.Cells.Item(110, 13).Value = -68433.002
$ws.Cells.Item(111, 8).Value = 2749.5
$ws.Cells.Item(111, 9).Value = 499
$ws.Cells.Item(111, 11).Value = 1497
$ws.Cells.Item(111, 13).Value = 1570
$ws.Cells.Item(114, 8).Value = 1147.3182
$ws.Cells.Item(114, 9).Value = 643.6429
$ws.Cells.Item(114, 10).Value = 2028.75
$ws.Cells.Item(114, 11).Value = 1930.9287
$ws.Cells.Item(114, 12).Value = 6086.25
$ws.Cells.Item(114, 13).Value = 1323.0713
$ws.Cells.Item(114, 14).Value = -12594.25
$ws.Cells.Item(115, 8).Value = 2271.4285
$ws.Cells.Item(115, 9).Value = 1725
$ws.Cells.Item(115, 11).Value = 5175
$ws.Cells.Item(115, 13).Value = -4000
$ws.Cells.Item(117, 8).Value = 1070.3684
$ws.Cells.Item(117, 9).Value = 752.2857
$ws.Cells.Item(117, 10).Value = 1255.9166
$ws.Cells.Item(117, 11).Value = 2256.8571
$ws.Cells.Item(117, 12).Value = 3767.7498
$ws.Cells.Item(117, 13).Value = 1185.1429
$ws.Cells.Item(117, 14).Value = -10651.7498
$ws.Cells.Item(118, 8).Value = 3833.1667
$ws.Cells.Item(118, 9).Value = 3799.8
$ws.Cells.Item(118, 11).Value = 11399.4
$ws.Cells.Item(118, 13).Value = -10156.4
$ws.Cells.Item(121, 8).Value = 604
$ws.Cells.Item(121, 10).Value = 1100
$ws.Cells.Item(121, 12).Value = 3300
$ws.Cells.Item(121, 14).Value = -5920
$ws.Cells.Item(127, 8).Value = 11812
$ws.Cells.Item(127, 10).Value = 11812
$ws.Cells.Item(127, 12).Value = 35436
$ws.Cells.Item(127, 14).Value = -45356
$ws.Cells.Item(133, 8).Value = 10226.25
$ws.Cells.Item(133, 9).Value = 9032.8125
$ws.Cells.Item(133, 11).Value = 27098.4375
$ws.Cells.Item(133, 13).Value = -22038.4375
$ws.Cells.Item(134, 8).Value = 7818.8823
$ws.Cells.Item(134, 9).Value = 5743.4165
$ws.Cells.Item(134, 10).Value = 12800
$ws.Cells.Item(134, 11).Value = 17230.2495
$ws.Cells.Item(134, 12).Value = 38400
$ws.Cells.Item(134, 13).Value = -12160.2495
$ws.Cells.Item(134, 14).Value = -48540
$ws.Cells.Item(135, 8).Value = 826.6429
$ws.Cells.Item(135, 9).Value = 826.6429
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 7439.7861
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 13).Value = -4904.7861
$ws.Cells.Item(135, 14).Value = ""
$ws.Cells.Item(136, 8).Value = 482988.53
$ws.Cells.Item(136, 9).Value = 717054.2
$ws.Cells.Item(136, 10).Value = 14857.143
$ws.Cells.Item(136, 11).Value = 2151162.6
$ws.Cells.Item(136, 12).Value = 44571.429
$ws.Cells.Item(136, 13).Value = -2146062.6
$ws.Cells.Item(136, 14).Value = -54771.429
$ws.Cells.Item(137, 8).Value = 1766.6
$ws.Cells.Item(137, 9).Value = 901.9
$ws.Cells.Item(137, 10).Value = 2631.3
$ws.Cells.Item(137, 11).Value = 2705.7
$ws.Cells.Item(137, 12).Value = 7893.900000000001
$ws.Cells.Item(137, 13).Value = 2394.3
$ws.Cells.Item(137, 14).Value = -18093.9
$ws.Cells.Item(138, 10).Value = 11499.25
$ws.Cells.Item(138, 12).Value = 34497.75
$ws.Cells.Item(138, 14).Value = -44777.75
$ws.Cells.Item(139, 8).Value = 2934
$ws.Cells.Item(139, 9).Value = 2925.75
$ws.Cells.Item(139, 11).Value = 8777.25
$ws.Cells.Item(139, 13).Value = -3637.25
$ws.Cells.Item(140, 8).Value = 1826.909
$ws.Cells.Item(140, 9).Value = 1826.909
$ws.Cells.Item(140, 11).Value = 5480.727000000001
$ws.Cells.Item(140, 13).Value = -300.7270000000008
$ws.Cells.Item(141, 8).Value = 55907.855
$ws.Cells.Item(141, 9).Value = 1005.3333
$ws.Cells.Item(141, 11).Value = 3015.9999
$ws.Cells.Item(141, 13).Value = 2164.0001
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(41, 8).Value = 3666.3333
$ws.Cells.Item(41, 10).Value = 0
$ws.Cells.Item(41, 12).Value = 0
$ws.Cells.Item(41, 14).Value = ""
$ws.Cells.Item(70, 8).Value = 20622.455
$ws.Cells.Item(70, 10).Value = 21262.125
$ws.Cells.Item(70, 12).Value = 21262.125
$ws.Cells.Item(70, 14).Value = -21802.125
$ws.Cells.Item(73, 8).Value = 20622.455
$ws.Cells.Item(73, 10).Value = 21262.125
$ws.Cells.Item(73, 12).Value = 21262.125
$ws.Cells.Item(73, 14).Value = -23134.125
$ws.Cells.Item(80, 8).Value = 2750
$ws.Cells.Item(80, 9).Value = 2750
$ws.Cells.Item(80, 10).Value = 2750
$ws.Cells.Item(80, 11).Value = 2750
$ws.Cells.Item(80, 12).Value = 2750
$ws.Cells.Item(80, 13).Value = -1752
$ws.Cells.Item(80, 14).Value = -4746
$ws.Cells.Item(83, 8).Value = 2750
$ws.Cells.Item(83, 9).Value = 2750
$ws.Cells.Item(83, 10).Value = 2750
$ws.Cells.Item(83, 11).Value = 13750
$ws.Cells.Item(83, 12).Value = 13750
$ws.Cells.Item(83, 13).Value = -8758
$ws.Cells.Item(83, 14).Value = -23734
$ws.Cells.Item(97, 8).Value = 1160.2222
$ws.Cells.Item(97, 9).Value = 1192.1666
$ws.Cells.Item(97, 10).Value = 1096.3334
$ws.Cells.Item(97, 11).Value = 1192.1666
$ws.Cells.Item(97, 12).Value = 1096.3334
$ws.Cells.Item(97, 13).Value = -696.1666
$ws.Cells.Item(97, 14).Value = -2088.3334
$ws.Cells.Item(122, 8).Value = 2985.853
$ws.Cells.Item(122, 9).Value = 2621.1
$ws.Cells.Item(122, 10).Value = 3506.9285
$ws.Cells.Item(122, 11).Value = 7863.299999999999
$ws.Cells.Item(122, 12).Value = 10520.7855
$ws.Cells.Item(122, 13).Value = -5413.299999999999
$ws.Cells.Item(122, 14).Value = -15420.7855
$ws.Cells.Item(126, 8).Value = 12635.954
$ws.Cells.Item(126, 9).Value = 16262
$ws.Cells.Item(126, 10).Value = 2966.5
$ws.Cells.Item(126, 11).Value = 48786
$ws.Cells.Item(126, 12).Value = 8899.5
$ws.Cells.Item(126, 13).Value = -46316
$ws.Cells.Item(126, 14).Value = -13839.5
$ws.Cells.Item(132, 8).Value = 4465.409
$ws.Cells.Item(132, 9).Value = 3577.5
$ws.Cells.Item(132, 10).Value = 6833.1665
$ws.Cells.Item(132, 11).Value = 10732.5
$ws.Cells.Item(132, 12).Value = 20499.4995
$ws.Cells.Item(132, 13).Value = -8202.5
$ws.Cells.Item(132, 14).Value = -25559.4995
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1600.9791
$ws.Cells.Item(16, 9).Value = 1236.8536
$ws.Cells.Item(16, 11).Value = 1236.8536
$ws.Cells.Item(16, 13).Value = -1066.8536
$ws.Cells.Item(22, 8).Value = 2313.3333
$ws.Cells.Item(22, 9).Value = 1694.3
$ws.Cells.Item(22, 10).Value = 2755.5
$ws.Cells.Item(22, 11).Value = 1694.3
$ws.Cells.Item(22, 12).Value = 2755.5
$ws.Cells.Item(22, 13).Value = -1399.3
$ws.Cells.Item(22, 14).Value = -3345.5
$ws.Cells.Item(27, 8).Value = 2313.3333
$ws.Cells.Item(27, 9).Value = 1694.3
$ws.Cells.Item(27, 10).Value = 2755.5
$ws.Cells.Item(27, 11).Value = 1694.3
$ws.Cells.Item(27, 12).Value = 2755.5
$ws.Cells.Item(27, 13).Value = -1587.3
$ws.Cells.Item(27, 14).Value = -2969.5
$ws.Cells.Item(40, 8).Value = 6749.1665
$ws.Cells.Item(40, 9).Value = 5998.3335
$ws.Cells.Item(40, 10).Value = 7500
$ws.Cells.Item(40, 11).Value = 5998.3335
$ws.Cells.Item(40, 12).Value = 7500
$ws.Cells.Item(40, 13).Value = -5862.3335
$ws.Cells.Item(40, 14).Value = -7772
$ws.Cells.Item(68, 8).Value = 2673.6428
$ws.Cells.Item(68, 10).Value = 1998.25
$ws.Cells.Item(68, 12).Value = 1998.25
$ws.Cells.Item(68, 14).Value = -3496.25
$ws.Cells.Item(71, 8).Value = 2673.6428
$ws.Cells.Item(71, 10).Value = 1998.25
$ws.Cells.Item(71, 12).Value = 9991.25
$ws.Cells.Item(71, 14).Value = -17479.25
$ws.Cells.Item(82, 8).Value = 5864.5454
$ws.Cells.Item(82, 9).Value = 6967.6665
$ws.Cells.Item(82, 10).Value = 900.5
$ws.Cells.Item(82, 11).Value = 6967.6665
$ws.Cells.Item(82, 12).Value = 900.5
$ws.Cells.Item(82, 13).Value = -6606.6665
$ws.Cells.Item(82, 14).Value = -1622.5
$ws.Cells.Item(85, 8).Value = 5864.5454
$ws.Cells.Item(85, 9).Value = 6967.6665
$ws.Cells.Item(85, 10).Value = 900.5
$ws.Cells.Item(85, 11).Value = 6967.6665
$ws.Cells.Item(85, 12).Value = 900.5
$ws.Cells.Item(85, 13).Value = -5719.6665
$ws.Cells.Item(85, 14).Value = -3396.5
$ws.Cells.Item(93, 8).Value = 1147.2858
$ws.Cells.Item(93, 9).Value = 1181.7142
$ws.Cells.Item(93, 10).Value = 1009.5714
$ws.Cells.Item(93, 11).Value = 1181.7142
$ws.Cells.Item(93, 12).Value = 1009.5714
$ws.Cells.Item(93, 13).Value = 66.28580000000011
$ws.Cells.Item(93, 14).Value = -3505.5714
$ws.Cells.Item(101, 8).Value = 18560.8
$ws.Cells.Item(101, 10).Value = 18560.8
$ws.Cells.Item(101, 12).Value = 18560.8
$ws.Cells.Item(101, 14).Value = -25050.8
$ws.Cells.Item(105, 8).Value = 27499.5
$ws.Cells.Item(105, 10).Value = 27499.5
$ws.Cells.Item(105, 12).Value = 27499.5
$ws.Cells.Item(105, 14).Value = -34487.5
$ws.Cells.Item(109, 8).Value = 54500
$ws.Cells.Item(109, 10).Value = 54500
$ws.Cells.Item(109, 12).Value = 54500
$ws.Cells.Item(109, 14).Value = -57274
$ws.Cells.Item(122, 8).Value = 559842.1
$ws.Cells.Item(122, 9).Value = 1004165.9
$ws.Cells.Item(122, 10).Value = 4437.375
$ws.Cells.Item(122, 11).Value = 3012497.7
$ws.Cells.Item(122, 12).Value = 13312.125
$ws.Cells.Item(122, 13).Value = -3010047.7
$ws.Cells.Item(122, 14).Value = -18212.125
$ws.Cells.Item(132, 8).Value = 4563.1816
$ws.Cells.Item(132, 9).Value = 4462.5
$ws.Cells.Item(132, 11).Value = 13387.5
$ws.Cells.Item(132, 13).Value = -10857.5
$ws.Cells.Item(135, 8).Value = 103000
$ws.Cells.Item(135, 10).Value = 103000
$ws.Cells.Item(135, 12).Value = 103000
$ws.Cells.Item(135, 14).Value = -113140
$ws.Cells.Item(136, 8).Value = 3190
$ws.Cells.Item(136, 9).Value = 3019.1667
$ws.Cells.Item(136, 11).Value = 9057.500100000001
$ws.Cells.Item(136, 13).Value = -6507.500100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 19648
$ws.Cells.Item(39, 9).Value = 18944
$ws.Cells.Item(39, 11).Value = 18944
$ws.Cells.Item(39, 13).Value = -18531
$ws.Cells.Item(40, 8).Value = 22200
$ws.Cells.Item(40, 10).Value = 22200
$ws.Cells.Item(40, 12).Value = 22200
$ws.Cells.Item(40, 14).Value = -22498
$ws.Cells.Item(75, 8).Value = 44000
$ws.Cells.Item(75, 10).Value = 40000
$ws.Cells.Item(75, 12).Value = 40000
$ws.Cells.Item(75, 14).Value = -41872
$ws.Cells.Item(78, 8).Value = 44000
$ws.Cells.Item(78, 10).Value = 40000
$ws.Cells.Item(78, 12).Value = 120000
$ws.Cells.Item(78, 14).Value = -129360
$ws.Cells.Item(105, 8).Value = 17358
$ws.Cells.Item(105, 10).Value = 17358
$ws.Cells.Item(105, 12).Value = 17358
$ws.Cells.Item(105, 14).Value = -24346
$ws.Cells.Item(107, 8).Value = 40923.96
$ws.Cells.Item(107, 9).Value = 954
$ws.Cells.Item(107, 11).Value = 2862
$ws.Cells.Item(107, 13).Value = -942
$ws.Cells.Item(113, 8).Value = 472.5
$ws.Cells.Item(113, 9).Value = 486.95834
$ws.Cells.Item(113, 10).Value = 429.125
$ws.Cells.Item(113, 11).Value = 1460.87502
$ws.Cells.Item(113, 12).Value = 1287.375
$ws.Cells.Item(113, 13).Value = 709.12498
$ws.Cells.Item(113, 14).Value = -5627.375
$ws.Cells.Item(122, 8).Value = 1843.4642
$ws.Cells.Item(122, 9).Value = 1793.48
$ws.Cells.Item(122, 10).Value = 2260
$ws.Cells.Item(122, 11).Value = 5380.440000000001
$ws.Cells.Item(122, 12).Value = 6780
$ws.Cells.Item(122, 13).Value = -2930.440000000001
$ws.Cells.Item(122, 14).Value = -11680
$ws.Cells.Item(123, 8).Value = 99999
$ws.Cells.Item(123, 10).Value = 99999
$ws.Cells.Item(123, 12).Value = 99999
$ws.Cells.Item(123, 14).Value = -109799
$ws.Cells.Item(126, 8).Value = 15424.556
$ws.Cells.Item(126, 9).Value = 16727.625
$ws.Cells.Item(126, 11).Value = 50182.875
$ws.Cells.Item(126, 13).Value = -47712.875
$ws.Cells.Item(131, 8).Value = 147332.67
$ws.Cells.Item(131, 10).Value = 147332.67
$ws.Cells.Item(131, 12).Value = 147332.67
$ws.Cells.Item(131, 14).Value = -157412.67
$ws.Cells.Item(132, 8).Value = 2163.3057
$ws.Cells.Item(132, 9).Value = 1976.8125
$ws.Cells.Item(132, 11).Value = 5930.4375
$ws.Cells.Item(132, 13).Value = -3400.4375
